$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy the number formatting from the (old D, now shifted to) E column
# onto the new D column so the new column matches the existing look
# (date format on the header rows, number format on the data rows).
# Done per data block so we don't create stray cells in the blank
# separator rows (6, 36, 37, 78, 79) that have no D/E content.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# Populate the new column D with the latest reporting period values
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 10553000
$ws.Range("D9").Value = 3505000
$ws.Range("D10").Value = 7048000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 1180000
$ws.Range("D15").Value = 1398000
$ws.Range("D17").Value = 8619000
$ws.Range("D18").Value = 1934000
$ws.Range("D20").Value = -183000
$ws.Range("D21").Value = 6437000
$ws.Range("D22").Value = 729000
$ws.Range("D23").Value = 1022000
$ws.Range("D24").Value = 341000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 681000
$ws.Range("D27").Value = 429000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 183000
$ws.Range("D33").Value = 429000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 429000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 986000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 2620000
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 625000
$ws.Range("D46").Value = 4231000
$ws.Range("D47").Value = 935000
$ws.Range("D48").Value = 800000
$ws.Range("D49").Value = 25749000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 835000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 32550000
$ws.Range("D57").Value = 325000
$ws.Range("D58").Value = 1860000
$ws.Range("D59").Value = 1812000
$ws.Range("D60").Value = 3997000
$ws.Range("D61").Value = 15185000
$ws.Range("D62").Value = 2851000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 24164000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 5254000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 8386000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 429000
$ws.Range("D83").Value = 4686000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 2576000
$ws.Range("D91").Value = -147000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -8593000
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -283000
$ws.Range("D101").Value = -23000
$ws.Range("D102").Value = -6323000
